# Update "想去人数" (want-to-go count) values in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new F-column value
$updates = @{
    2  = 622
    3  = 575
    6  = 116
    8  = 57
    10 = 4998
    11 = 4697
    12 = 14
    15 = 41
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
